# Auto-generated script to update computed market-price / profit columns
# on several sheets, matching a scheduled data refresh from the Universalis API.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 877.6667
$ws.Range("I28").Value = 966.5294
$ws.Range("K28").Value = 966.5294
$ws.Range("M28").Value = -481.5294
$ws.Range("H106").Value = 6608.5
$ws.Range("I106").Value = 1256.0769
$ws.Range("K106").Value = 1256.0769
$ws.Range("M106").Value = -625.0769
$ws.Range("H111").Value = 353.2857
$ws.Range("I111").Value = 353.2857
$ws.Range("K111").Value = 1059.8571
$ws.Range("M111").Value = 2007.1429
$ws.Range("H125").Value = 1479.3077
$ws.Range("I125").Value = 1313.2
$ws.Range("J125").Value = 2033
$ws.Range("K125").Value = 11818.8
$ws.Range("L125").Value = 18297
$ws.Range("M125").Value = -9358.800000000001
$ws.Range("N125").Value = -23217
$ws.Range("H137").Value = 1329.4667
$ws.Range("I137").Value = 1276.3077
$ws.Range("J137").Value = 1675
$ws.Range("K137").Value = 3828.9231
$ws.Range("L137").Value = 5025
$ws.Range("M137").Value = -1278.9231
$ws.Range("N137").Value = -10125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1510.4762
$ws.Range("I2").Value = 1505.25
$ws.Range("J2").Value = 1541.8334
$ws.Range("K2").Value = 1505.25
$ws.Range("L2").Value = 1541.8334
$ws.Range("M2").Value = -1392.25
$ws.Range("N2").Value = -1767.8334
$ws.Range("H32").Value = 6184.92
$ws.Range("I32").Value = 5250
$ws.Range("K32").Value = 5250
$ws.Range("M32").Value = -4963
$ws.Range("H116").Value = 1510.4762
$ws.Range("I116").Value = 1505.25
$ws.Range("J116").Value = 1541.8334
$ws.Range("K116").Value = 1505.25
$ws.Range("L116").Value = 1541.8334
$ws.Range("M116").Value = 788.75
$ws.Range("N116").Value = -6129.8334
$ws.Range("H122").Value = 1829.3864
$ws.Range("I122").Value = 1429.5862
$ws.Range("J122").Value = 2602.3333
$ws.Range("K122").Value = 4288.7586
$ws.Range("L122").Value = 7806.999899999999
$ws.Range("M122").Value = -1838.7586
$ws.Range("N122").Value = -12706.9999
$ws.Range("H126").Value = 4569.533
$ws.Range("I126").Value = 4569.533
$ws.Range("K126").Value = 13708.599
$ws.Range("M126").Value = -11238.599
$ws.Range("H132").Value = 1540.356
$ws.Range("I132").Value = 1086.6852
$ws.Range("J132").Value = 6440
$ws.Range("K132").Value = 3260.0556
$ws.Range("L132").Value = 19320
$ws.Range("M132").Value = -730.0555999999997
$ws.Range("N132").Value = -24380

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1510.4762
$ws.Range("I3").Value = 1505.25
$ws.Range("J3").Value = 1541.8334
$ws.Range("K3").Value = 1505.25
$ws.Range("L3").Value = 1541.8334
$ws.Range("M3").Value = -1391.25
$ws.Range("N3").Value = -1769.8334
$ws.Range("H107").Value = 70085.92999999999
$ws.Range("I107").Value = 112088.22
$ws.Range("J107").Value = 7082.5
$ws.Range("K107").Value = 112088.22
$ws.Range("L107").Value = 7082.5
$ws.Range("M107").Value = -110168.22
$ws.Range("N107").Value = -10922.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 71428860
$ws.Range("I7").Value = 200000220
$ws.Range("J7").Value = 327.55554
$ws.Range("K7").Value = 200000220
$ws.Range("L7").Value = 327.55554
$ws.Range("M7").Value = -200000107
$ws.Range("N7").Value = -553.5555400000001
$ws.Range("H107").Value = 1915.069
$ws.Range("I107").Value = 1534.9445
$ws.Range("K107").Value = 1534.9445
$ws.Range("M107").Value = 385.0554999999999
$ws.Range("H132").Value = 3449.8
$ws.Range("I132").Value = 3966.0476
$ws.Range("J132").Value = 2245.2222
$ws.Range("K132").Value = 11898.1428
$ws.Range("L132").Value = 6735.6666
$ws.Range("M132").Value = -9368.1428
$ws.Range("N132").Value = -11795.6666
$ws.Range("H134").Value = 933.2558
$ws.Range("I134").Value = 933.2558
$ws.Range("K134").Value = 2799.7674
$ws.Range("M134").Value = -264.7674000000002
$ws.Range("H138").Value = 108998
$ws.Range("J138").Value = 112797.6
$ws.Range("L138").Value = 112797.6
$ws.Range("N138").Value = -123077.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9864.125
$ws.Range("I56").Value = 9864.125
$ws.Range("K56").Value = 9864.125
$ws.Range("M56").Value = -9334.125
$ws.Range("H113").Value = 2302.5625
$ws.Range("J113").Value = 2302.5625
$ws.Range("L113").Value = 6907.6875
$ws.Range("N113").Value = -11247.6875
$ws.Range("H131").Value = 2847.1875
$ws.Range("I131").Value = 1694.875
$ws.Range("K131").Value = 5084.625
$ws.Range("M131").Value = -44.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 9091384
$ws.Range("J3").Value = 11429282
$ws.Range("L3").Value = 11429282
$ws.Range("N3").Value = -11429514
$ws.Range("H12").Value = 516.6667
$ws.Range("J12").Value = 516.6667
$ws.Range("L12").Value = 516.6667
$ws.Range("N12").Value = -796.6667
$ws.Range("H113").Value = 3836.1177
$ws.Range("J113").Value = 5772.4
$ws.Range("L113").Value = 5772.4
$ws.Range("N113").Value = -10112.4
$ws.Range("H132").Value = 1881.0416
$ws.Range("I132").Value = 1793.6818
$ws.Range("J132").Value = 2842
$ws.Range("K132").Value = 5381.0454
$ws.Range("L132").Value = 8526
$ws.Range("M132").Value = -2851.0454
$ws.Range("N132").Value = -13586
$ws.Range("H133").Value = 94839
$ws.Range("J133").Value = 94839
$ws.Range("L133").Value = 94839
$ws.Range("N133").Value = -104959
$ws.Range("H136").Value = 30509.352
$ws.Range("J136").Value = 30509.352
$ws.Range("L136").Value = 91528.056
$ws.Range("N136").Value = -96628.056

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 4250
$ws.Range("J11").Value = 4250
$ws.Range("L11").Value = 4250
$ws.Range("N11").Value = -4530
$ws.Range("H16").Value = 1193.1333
$ws.Range("I16").Value = 1138.3077
$ws.Range("K16").Value = 1138.3077
$ws.Range("M16").Value = -968.3077000000001
$ws.Range("H46").Value = 34341.77
$ws.Range("I46").Value = 61874
$ws.Range("K46").Value = 61874
$ws.Range("M46").Value = -61686
$ws.Range("H132").Value = 2642.7207
$ws.Range("I132").Value = 2393.4905
$ws.Range("J132").Value = 3523.3333
$ws.Range("K132").Value = 7180.4715
$ws.Range("L132").Value = 10569.9999
$ws.Range("M132").Value = -4650.4715
$ws.Range("N132").Value = -15629.9999
$ws.Range("H136").Value = 3926.0833
$ws.Range("I136").Value = 3680.889
$ws.Range("J136").Value = 4661.6665
$ws.Range("K136").Value = 11042.667
$ws.Range("L136").Value = 13984.9995
$ws.Range("M136").Value = -8492.667000000001
$ws.Range("N136").Value = -19084.9995
$ws.Range("H141").Value = 105357
$ws.Range("J141").Value = 105357
$ws.Range("L141").Value = 105357
$ws.Range("N141").Value = -115717

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 4086
$ws.Range("I9").Value = 5379.25
$ws.Range("J9").Value = 1499.5
$ws.Range("K9").Value = 5379.25
$ws.Range("L9").Value = 1499.5
$ws.Range("M9").Value = -5239.25
$ws.Range("N9").Value = -1779.5
$ws.Range("H107").Value = 16131529
$ws.Range("I107").Value = 1802.6842
$ws.Range("J107").Value = 41670260
$ws.Range("K107").Value = 5408.0526
$ws.Range("L107").Value = 125010780
$ws.Range("M107").Value = -3488.0526
$ws.Range("N107").Value = -125014620
$ws.Range("H132").Value = 2719.2778
$ws.Range("I132").Value = 2159.6123
$ws.Range("J132").Value = 8204
$ws.Range("K132").Value = 6478.836899999999
$ws.Range("L132").Value = 24612
$ws.Range("M132").Value = -3948.836899999999
$ws.Range("N132").Value = -29672
$ws.Range("H136").Value = 917.73334
$ws.Range("I136").Value = 989.5
$ws.Range("K136").Value = 2968.5
$ws.Range("M136").Value = -418.5
